# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") previously held a raw strikeout-total style figure; it is
# regenerated here to the corrected per-game K value for each of the 25
# logged outings (rows 2-26 on the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 6
    4  = 6
    5  = 2
    6  = 8
    7  = 7
    8  = 7
    9  = 3
    10 = 9
    11 = 5
    12 = 8
    13 = 6
    14 = 6
    15 = 9
    16 = 5
    17 = 5
    18 = 10
    19 = 5
    20 = 0
    21 = 4
    22 = 6
    23 = 6
    24 = 4
    25 = 4
    26 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
